# Edit script: add culturalgroup_male / culturalgroup_female sheets,
# populate data, apply custom font style to column B, fix up tab/view state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the two new worksheets at the end of the tab strip, in order:
#    ... choices_male, culturalgroup_male, culturalgroup_female
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMale = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsMale.Name = "culturalgroup_male"
$wsFemale = $wb.Worksheets.Add([Type]::Missing, $wsMale)
$wsFemale.Name = "culturalgroup_female"

# ---------------------------------------------------------------------------
# 2. Register a one-off style carrying the "BitstromWera Nerd Font" font so
#    the engine allocates a single clean font entry (matches styles.xml diff:
#    fonts count 2 -> 3, cellXfs count 1 -> 2). We delete the named style
#    again afterwards so cellStyles / cellStyleXfs stay untouched (count=1),
#    exactly like the target - only the xf + font survive.
# ---------------------------------------------------------------------------
$culturalFontStyle = $wb.Styles.Add("CulturalGroupFont")
$culturalFontStyle.Font.Name = "BitstromWera Nerd Font"

# ---------------------------------------------------------------------------
# 3. culturalgroup_male data (36 rows incl. header)
# ---------------------------------------------------------------------------
$maleData = @(
    @("student", "cultural group"),
    @("A", "ASIA"),
    @("B", "ASIA"),
    @("C", "ASIA"),
    @("D", "ESEU"),
    @("E", "ESEU"),
    @("F", "ESEU"),
    @("G", "MENA"),
    @("H", "MENA"),
    @("I", "MENA"),
    @("J", "MENA"),
    @("K", "LTIN"),
    @("L", "LTIN"),
    @("M", "ITAL"),
    @("N", "ITAL"),
    @("O", "ITAL"),
    @("P", "ITAL"),
    @("Q", "ITAL"),
    @("R", "ITAL"),
    @("S", "ITAL"),
    @("T", "ITAL"),
    @("U", "ALTR"),
    @("V", "ALTR"),
    @("W", "ALTR"),
    @("X", "ALTR"),
    @("Y", "ALTR"),
    @("Z", "ALTR"),
    @("AA", "ALTR"),
    @("AB", "ALTR"),
    @("AC", "ALTR"),
    @("AD", "ALTR"),
    @("AE", "ALTR"),
    @("AF", "ALTR"),
    @("AG", "ALTR"),
    @("AH", "ALTR"),
    @("AI", "ALTR")
)

for ($i = 0; $i -lt $maleData.Length; $i++) {
    $r = $i + 1
    $wsMale.Cells.Item($r, 1).Value = $maleData[$i][0]
    $wsMale.Cells.Item($r, 2).Value = $maleData[$i][1]
    $wsMale.Cells.Item($r, 2).Style = "CulturalGroupFont"
}
# header row's A cell has no explicit style in the diff (A1..A36 unstyled),
# only column B carries the custom font.

# trailing styled-but-empty B cells, rows 37-57
for ($r = 37; $r -le 57; $r++) {
    $wsMale.Cells.Item($r, 2).Style = "CulturalGroupFont"
}

$wsMale.Columns.Item(2).ColumnWidth = 9.14

# ---------------------------------------------------------------------------
# 4. culturalgroup_female data (57 rows incl. header)
# ---------------------------------------------------------------------------
$femaleData = @(
    @("student", "cultural group"),
    @("Female1", "ASIA"),
    @("Female2", "ASIA"),
    @("Female3", "ASIA"),
    @("Female4", "ASIA"),
    @("Female5", "ASIA"),
    @("Female6", "ASIA"),
    @("Female7", "ASIA"),
    @("Female8", "ASIA"),
    @("Female9", "ASIA"),
    @("Female10", "ASIA"),
    @("Female11", "ASIA"),
    @("Female12", "ESEU"),
    @("Female13", "ESEU"),
    @("Female14", "ESEU"),
    @("Female15", "ESEU"),
    @("Female16", "ESEU"),
    @("Female17", "ESEU"),
    @("Female18", "MENA"),
    @("Female19", "MENA"),
    @("Female20", "MENA"),
    @("Female21", "MENA"),
    @("Female22", "LTIN"),
    @("Female23", "LTIN"),
    @("Female24", "ITAL"),
    @("Female25", "ITAL"),
    @("Female26", "ITAL"),
    @("Female27", "ITAL"),
    @("Female28", "ITAL"),
    @("Female29", "ITAL"),
    @("Female30", "ITAL"),
    @("Female31", "ITAL"),
    @("Female32", "ALTR"),
    @("Female33", "ALTR"),
    @("Female34", "ALTR"),
    @("Female35", "ALTR"),
    @("Female36", "ALTR"),
    @("Female37", "ALTR"),
    @("Female38", "ALTR"),
    @("Female39", "ALTR"),
    @("Female40", "ALTR"),
    @("Female41", "ALTR"),
    @("Female42", "ALTR"),
    @("Female43", "ALTR"),
    @("Female44", "ALTR"),
    @("Female45", "ALTR"),
    @("Female46", "ALTR"),
    @("Female47", "ALTR"),
    @("Female48", "ALTR"),
    @("Female49", "ALTR"),
    @("Female50", "ALTR"),
    @("Female51", "ALTR"),
    @("Female52", "ALTR"),
    @("Female53", "ALTR"),
    @("Female54", "ALTR"),
    @("Female55", "ALTR"),
    @("Female56", "ALTR")
)

for ($i = 0; $i -lt $femaleData.Length; $i++) {
    $r = $i + 1
    $wsFemale.Cells.Item($r, 1).Value = $femaleData[$i][0]
    $wsFemale.Cells.Item($r, 1).Style = "CulturalGroupFont"
    $wsFemale.Cells.Item($r, 2).Value = $femaleData[$i][1]
    $wsFemale.Cells.Item($r, 2).Style = "CulturalGroupFont"
}

$wsFemale.Columns.Item(1).ColumnWidth = 13
$wsFemale.Columns.Item(2).ColumnWidth = 17.71

# Drop the helper named style now that every cell referencing it has been
# written - leaves a bare font + cellXf behind, cellStyles/cellStyleXfs
# collapse back to just "Normal" (matches the target styles.xml exactly).
$culturalFontStyle.Delete()

# ---------------------------------------------------------------------------
# 5. View-state touch-ups on the pre-existing sheets
# ---------------------------------------------------------------------------
$wsCapacityFemale = $wb.Worksheets.Item("capacity_female")
$wsCapacityFemale.Activate()
$wsCapacityFemale.Range("D11").Select()

$wsChoicesMale = $wb.Worksheets.Item("choices_male")
$wsChoicesMale.Activate()
$wsChoicesMale.Range("B3").Select()

# culturalgroup_female: scrolled near the bottom, selection sits on B1:B57
$wsFemale.Activate()
$excel.ActiveWindow.ScrollRow = 37
$wsFemale.Range("B1:B57").Select()

# culturalgroup_male ends up the active tab, scrolled down a bit, C20 selected
$wsMale.Activate()
$excel.ActiveWindow.ScrollRow = 8
$wsMale.Range("C20").Select()

# first visible tab in the tab strip is choices_female (index 3, 0-based 2)
$excel.ActiveWindow.DisplayedSheets = $wb.Worksheets.Item("choices_female")
